$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits right after the italic
#    "Cardgame" run in the Overview paragraph. In the edited
#    document it is relocated to just before the final "." of the
#    brand-new "Players" paragraph being added below, so drop it
#    from its old spot first (it is recreated via the inserted XML).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Locate the "Players" heading paragraph.
# ------------------------------------------------------------------
$playersHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Players") {
        $playersHeadingIndex = $i
        break
    }
}

$playersPara = $d.Paragraphs.Item($playersHeadingIndex)

# ------------------------------------------------------------------
# 3. Open up a brand-new (empty, style-less) paragraph right after
#    the heading, then pour in the authored WordOpenXML so it ends
#    up with exactly the run/formatting/bookmark layout of the real
#    edit (no style inherited from the heading, proofErr spell-check
#    markers kept around "Cardgame", and the relocated "_GoBack"
#    bookmark positioned right before the trailing period).
# ------------------------------------------------------------------
$playersPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($playersHeadingIndex + 1)
$newPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>Cardgame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is a game for two players who go head-to-head with decks that they’ve created. Each card in the deck represents either a unit in their army (for example, a plague doctor) or a utility that they can use to get ahead (for example, a bank withdrawal). The exception to this is the General. Each player can select one of a number of General’s, which acts as not only their character with a unique ability, but also as their avatar. This would be similar to the class a player selects in </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>Hearthstone</w:t></w:r><w:r><w:t xml:space="preserve">, or their Commander in the Commander format of </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>Magic: The Gathering</w:t></w:r><w:r><w:t xml:space="preserve"> (however without the limitation to deck-building present in both of these cases).</w:t></w:r><w:r><w:t xml:space="preserve"> Player’s would play together over a network on two separate machines, however the game is designed in such a way that it would be easily replicable in real life with physical cards. To specify, in-game the player can only do with the cards what they’d be able to do with a card in real life, for example, change its</w:t></w:r><w:r><w:t xml:space="preserve"> orientation,</w:t></w:r><w:r><w:t xml:space="preserve"> flip it over</w:t></w:r><w:r><w:t>, or put counters on it</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

Write-Output "done"
